$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "(something) ends (～が)"
$ws.Range("A6").Value = "(something) begins (～が)"
$ws.Range("A7").Value = "to play (a string instrument or piano) (～を)"
$ws.Range("A8").Value = "to get (from somebody) (person に thing を)"
$ws.Range("A9").Value = "to memorize (～を)"
$ws.Range("A10").Value = "to appear; to attend (～に); to exit (～を)"
